{"js": "// Apply the \"Play Free Elephant King...\" -> new title/tagline rewrite,\n// plus the refreshed pros/cons bullet list text.\n// Each entry: [old text to find, new text to insert]\nconst replacements = [\n  [\n    \"Play Free Elephant King Online Slot Game Review\",\n    \"Play Elephant King Free - Elegant African-themed Slot Game\",\n  ],\n  [\n    \"Elegant graphics and sound effects in an African savanna\",\n    \"Elegant and refined graphics\",\n  ],\n  [\n    \"40 winning lines available for higher chance of hitting payout\",\n    \"Classic and innovative gameplay elements\",\n  ],\n  [\n    \"Chance to win random payouts adds an enticing element\",\n    \"Chance to win random payouts\",\n  ],\n  [\n    \"Combination of classic and innovative gameplay elements\",\n    \"Beautiful African-themed visuals and soundtrack\",\n  ],\n  [\n    \"Betting values can easily exceed player's budget\",\n    \"Betting values can quickly become high\",\n  ],\n  [\n    \"Free spins mode challenging to achieve\",\n    \"Achieving the free spins mode can be challenging\",\n  ],\n  [\n    \"Explore the beauty of African savanna with the Elephant King online slot game. Play for free and experience high potential for payouts with up to 40 winning lines.\",\n    \"Experience the elegance of Elephant King, a free online slot game with an African theme.\",\n  ],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Rewrite the title/tagline and refresh the pros/cons bullet list text.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old = \"Play Free Elephant King Online Slot Game Review\"; New = \"Play Elephant King Free - Elegant African-themed Slot Game\"},\n    @{Old = \"Elegant graphics and sound effects in an African savanna\"; New = \"Elegant and refined graphics\"},\n    @{Old = \"40 winning lines available for higher chance of hitting payout\"; New = \"Classic and innovative gameplay elements\"},\n    @{Old = \"Chance to win random payouts adds an enticing element\"; New = \"Chance to win random payouts\"},\n    @{Old = \"Combination of classic and innovative gameplay elements\"; New = \"Beautiful African-themed visuals and soundtrack\"},\n    @{Old = \"Betting values can easily exceed player's budget\"; New = \"Betting values can quickly become high\"},\n    @{Old = \"Free spins mode challenging to achieve\"; New = \"Achieving the free spins mode can be challenging\"},\n    @{Old = \"Explore the beauty of African savanna with the Elephant King online slot game. Play for free and experience high potential for payouts with up to 40 winning lines.\"; New = \"Experience the elegance of Elephant King, a free online slot game with an African theme.\"}\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
